$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Certificate")
$ws.Activate()

$ws.Range("B5").Value = "NAFA"
$ws.Range("A5").Value = "Painting"

$ws.Range("B5").Select()
